# Update the Ukrainian "Акт звіряння" (reconciliation act) blank to the
# Russian "Остатки на складе" (stock validation) blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Rename the worksheet/tab.
$ws.Name = "Остатки на складе"

# 2. Title row (merged A1:E1).
$ws.Range("A1").Value = "Остатки на складе ({`$v->cat_name}) на {`$v->date}"

# 3. "Залишок" -> "Остаток" header (merged D3:E3).
$ws.Range("D3").Value = "Остаток"

# 4. Signature lines.
$ws.Range("B8").Value = "Виписал________________________________ /{`$v->user_sign}/  "
$ws.Range("B9").Value = "Проверил_______________________________"

# 5. Move the active selection the way the authored file has it.
$ws.Range("D12").Select()
